# "Loan RBI, Variable Instalments"
#
# The "Repayment Schedule" sheet gets a new (blank) column inserted right
# before the old "Late" column (which, together with the old "Outstanding"
# column after it, shifts one slot to the right: N -> O, O -> P(unused),
# P -> Q). Column N1 becomes an empty/blank header cell; the old N-column
# "Late" values (all zero) land in column O, and the old P-column
# "Outstanding" values land in column Q.
#
# The workbook was also re-saved with the "Repayment Schedule" tab active
# and cell J22 selected there (previously "NewLoanInput" was the active
# tab, and F24 was selected on "Repayment Schedule").

$wb = $excel.ActiveWorkbook

$repaymentSchedule = $wb.Worksheets.Item("Repayment Schedule")

# Insert a new blank column at "N" (14th column), pushing the existing
# "Late" (N) and "Outstanding" (P) columns one slot to the right.
$repaymentSchedule.Columns.Item(14).Insert()

# Make "Repayment Schedule" the active/selected tab, with J22 selected -
# matching the workbook's new bookViews/activeTab + sheetView state.
$repaymentSchedule.Activate()
$repaymentSchedule.Range("J22").Select()
